$d = $word.ActiveDocument

$replacements = @(
    @("153×7=1071", "126×7=882"),
    @("437×3=1311", "889×7=6223"),
    @("766×4=3064", "360×9=3240"),
    @("255×7=1785", "602×7=4214"),
    @("403×3=1209", "969×4=3876"),
    @("425×5=2125", "567×5=2835"),
    @("903×4=3612", "685×5=3425"),
    @("413×7=2891", "177×9=1593"),
    @("955×7=6685", "406×3=1218"),
    @("135×6=810",  "702×9=6318"),
    @("879×7=6153", "531×4=2124"),
    @("406×4=1624", "991×2=1982"),
    @("101×2=202",  "241×4=964"),
    @("131×2=262",  "735×3=2205"),
    @("296×9=2664", "585×5=2925"),
    @("972×9=8748", "357×9=3213"),
    @("642×5=3210", "451×7=3157"),
    @("921×7=6447", "832×9=7488"),
    @("289×3=867",  "525×6=3150"),
    @("746×3=2238", "583×3=1749"),
    @("682×5=3410", "338×2=676"),
    @("498×5=2490", "151×9=1359"),
    @("534×9=4806", "176×5=880"),
    @("179×3=537",  "912×4=3648"),
    @("706×4=2824", "483×2=966")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
